# Refresh the cryptos list snapshot (Price / Volume(1h) columns) to the
# latest scraped values. Numeric-looking price strings are written with a
# leading apostrophe so Excel keeps them as text (matching the original
# inlineStr cell type) instead of silently coercing them to doubles, which
# would corrupt trailing zeros (e.g. "149.30" -> 149.3) and introduce
# binary floating point noise (e.g. "539.19" -> 539.19000000000005).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "59.391.27"
$ws.Range('E2').Value = "  +0.15%  "

$ws.Range('D3').Value = "2.604.67"
$ws.Range('E3').Value = "  +0.35%  "

$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = "  -0.12%  "

$ws.Range('D5').Value = "'539.19"
$ws.Range('E5').Value = "  +3.33%  "

$ws.Range('D6').Value = "'140.85"
$ws.Range('E6').Value = "  +0.37%  "

$ws.Range('E7').Value = "  +0.33%  "

$ws.Range('E8').Value = "  +0.29%  "

$ws.Range('E9').Value = "  -1.04%  "

$ws.Range('E10').Value = "  +1.63%  "

$ws.Range('E11').Value = "  +1.01%  "

$ws.Range('E12').Value = "  +1.88%  "

$ws.Range('D13').Value = "3.065.76"
$ws.Range('E13').Value = "  +0.39%  "

$ws.Range('D14').Value = "59.285.51"
$ws.Range('E14').Value = "  +0.25%  "

$ws.Range('D15').Value = "'20.54"
$ws.Range('E15').Value = "  +0.60%  "

$ws.Range('E16').Value = "  +0.73%  "

$ws.Range('D17').Value = "2.599.40"
$ws.Range('E17').Value = "  -0.32%  "

$ws.Range('D18').Value = "'343.07"
$ws.Range('E18').Value = "  +1.45%  "

$ws.Range('E19').Value = "  +0.47%  "

$ws.Range('D20').Value = "'10.11"
$ws.Range('E20').Value = "  -0.85%  "

$ws.Range('D21').Value = "'6.41"
$ws.Range('E21').Value = "  -1.15%  "

$ws.Range('E22').Value = "  +0.40%  "

$ws.Range('E23').Value = "  +1.91%  "

$ws.Range('E24').Value = "  -0.64%  "

$ws.Range('E25').Value = "  +1.08%  "

$ws.Range('E26').Value = "  +0.34%  "

$ws.Range('D27').Value = "'7.19"
$ws.Range('E27').Value = "  +1.92%  "

$ws.Range('E28').Value = "  +0.17%  "

$ws.Range('E29').Value = "  +1.68%  "

$ws.Range('E30').Value = "  +5.72%  "

$ws.Range('E31').Value = "  -1.97%  "

$ws.Range('E32').Value = "  -0.20%  "

$ws.Range('D33').Value = "'149.30"
$ws.Range('E33').Value = "  +0.09%  "

$ws.Range('E34').Value = "  -0.69%  "

$ws.Range('E35').Value = "  -1.12%  "

$ws.Range('D36').Value = "'36.97"
$ws.Range('E36').Value = "  +1.72%  "

$ws.Range('E37').Value = "  +0.39%  "

$ws.Range('E38').Value = "  +0.35%  "

$ws.Range('D39').Value = "'0.824"
$ws.Range('E39').Value = "  -0.06%  "

$ws.Range('E40').Value = "  +0.37%  "

$ws.Range('E41').Value = "  +0.45%  "

$ws.Range('D42').Value = "'273.80"
$ws.Range('E42').Value = "  -0.68%  "

$ws.Range('D43').Value = "'10.76"
$ws.Range('E43').Value = "  +0.35%  "

$ws.Range('E44').Value = "  +0.76%  "

$ws.Range('D45').Value = "'0.0958"
$ws.Range('E45').Value = "  +0.56%  "

$ws.Range('D46').Value = "'0.0523"
$ws.Range('E46').Value = "  +0.53%  "

$ws.Range('D47').Value = "1.946.71"

$ws.Range('E48').Value = "  +0.98%  "

$ws.Range('D49').Value = "'18.27"
$ws.Range('E49').Value = "  +0.55%  "

$ws.Range('E50').Value = "  -2.17%  "

$ws.Range('D51').Value = "'110.91"
$ws.Range('E51').Value = "  -2.01%  "
